# Update the "Out of PO" player list:
#  - Reorders the existing players
#  - Removes "Toumani Camara" and "Luguentz Dort"
# Resulting table occupies A1:C17 (16 data rows + header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Jonas Valanciunas", "C", "Washington Wizards"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Deandre Ayton", "C", "Portland Trail Blazers")
)

# Clear out the old data range (rows 2 through 19) before writing the new,
# shorter table so no stale rows (18-19) remain.
$oldDataRange = $ws.Range("A2:C19")
$oldDataRange.Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
